# "gps_alt_m changed from uint8_t to uint16_t" -- the Lipo 3S measurement
# series (H11:H13, "Digit") switches from an 8-bit ADC reading to a 16-bit
# one, so the raw digit counts and the derived "Lipo Spannung" (F column)
# values are updated for rows 11-13 on "Tabelle1". F11:F13 additionally
# pick up a new "0.000" number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New number format for the measured Lipo-cell voltages.
$ws.Range("F11:F13").NumberFormat = "0.000"

# Row 11 (cell 1)
$ws.Range("F11").Value = 4.1669999999999998
$ws.Range("H11").Value = 828

# Row 12 (cell 1+2)
$ws.Range("F12").Value = 4.1740000000000004
$ws.Range("H12").Value = 894

# Row 13 (cell 1+2+3)
$ws.Range("F13").Value = 4.1710000000000003
$ws.Range("H13").Value = 946

# Widen column F to fit the new values (it held a narrow, header-only
# width before).
$ws.Columns("F").AutoFit() | Out-Null

# Leave the selection where the user ended up after editing the table.
$ws.Range("J11:J13").Select()
